$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 3.6
$ws.Range("G2").Value = 4.3
$ws.Range("H2").Value = 2.08
$ws.Range("I2").Value = 2.32
$ws.Range("J2").Value = 3.2
$ws.Range("K2").Value = 3.75
$ws.Range("N2").Value = 2.36
$ws.Range("P2").Value = 1.63
$ws.Range("R2").Value = 1.2
$ws.Range("S2").Value = 2.68
$ws.Range("T2").Value = 1.04
$ws.Range("U2").Value = 1.04
$ws.Range("V2").Value = 1.75
$ws.Range("W2").Value = 1.3

# Row 3
$ws.Range("L3").Value = 1.36
$ws.Range("T3").Value = 1.85

# Row 4
$ws.Range("F4").Value = 2.66
$ws.Range("H4").Value = 2.42
$ws.Range("J4").Value = 2.66
$ws.Range("K4").Value = 4.7
$ws.Range("L4").Value = 1.35
$ws.Range("N4").Value = 2.52
$ws.Range("S4").Value = 2.72
$ws.Range("T4").Value = 1.04
$ws.Range("U4").Value = 1.04

# Row 5
$ws.Range("F5").Value = 1.49
$ws.Range("G5").Value = 1.63
$ws.Range("H5").Value = 6
$ws.Range("I5").Value = 9.6
$ws.Range("J5").Value = 4.1
$ws.Range("K5").Value = 5.6
$ws.Range("M5").Value = 1.05
$ws.Range("N5").Value = 3.7
$ws.Range("O5").Value = 1.27
$ws.Range("P5").Value = 1.93
$ws.Range("Q5").Value = 1.8
$ws.Range("R5").Value = 1.36
$ws.Range("S5").Value = 2.86
$ws.Range("T5").Value = 1.92
$ws.Range("U5").Value = 1.83
$ws.Range("V5").Value = 1.11
$ws.Range("W5").Value = 2.58
$ws.Range("AB5").Value = 9.6
$ws.Range("AF5").Value = 11
$ws.Range("AG5").Value = 12
$ws.Range("AN5").Value = 10.5

# Row 6
$ws.Range("I6").Value = 2.54
$ws.Range("L6").Value = 1.01
$ws.Range("N6").Value = 2.84
$ws.Range("S6").Value = 2.98
$ws.Range("T6").Value = 1.59
$ws.Range("U6").Value = 1.76

# Row 7
$ws.Range("G7").Value = 1.51
$ws.Range("L7").Value = 1.01
$ws.Range("N7").Value = 3.4
$ws.Range("R7").Value = 1.43
$ws.Range("S7").Value = 2.36
$ws.Range("T7").Value = 1.73
$ws.Range("U7").Value = 1.04
$ws.Range("V7").Value = 1.12
$ws.Range("W7").Value = 2.96

# Row 8
$ws.Range("F8").Value = 1.64
$ws.Range("G8").Value = 1.91
$ws.Range("H8").Value = 5
$ws.Range("I8").Value = 10.5
$ws.Range("J8").Value = 3.05
$ws.Range("K8").Value = 4.2
$ws.Range("L8").Value = 1.01
$ws.Range("M8").Value = 1.01
$ws.Range("N8").Value = 1.53
$ws.Range("O8").Value = 1.01
$ws.Range("P8").Value = 1.53
$ws.Range("Q8").Value = 2.3
$ws.Range("R8").Value = 1.14
$ws.Range("S8").Value = 3.9
$ws.Range("T8").Value = 1.94
$ws.Range("U8").Value = 1.43
$ws.Range("V8").Value = 1.12
$ws.Range("W8").Value = 2.08
$ws.Range("X8").Value = 13
$ws.Range("Y8").Value = 24
$ws.Range("Z8").Value = 1000
$ws.Range("AA8").Value = 1000
$ws.Range("AB8").Value = 8.199999999999999
$ws.Range("AC8").Value = 12.5
$ws.Range("AD8").Value = 1000
$ws.Range("AE8").Value = 1000
$ws.Range("AF8").Value = 12.5
$ws.Range("AG8").Value = 16
$ws.Range("AH8").Value = 1000
$ws.Range("AI8").Value = 1000
$ws.Range("AJ8").Value = 26
$ws.Range("AK8").Value = 1000
$ws.Range("AL8").Value = 1000
$ws.Range("AM8").Value = 1000
$ws.Range("AN8").Value = 1000
$ws.Range("AO8").Value = 1000

$wb.Save()